# Rename the element/material sector header labels (row 3, columns D-G) on
# every yearly sheet so they align with the Baseline naming:
#   Nd -> Neodymium
#   Dy -> Dysprosium
#   Cu -> Copper ores and concentrates
#   Si -> Raw silicon
#
# These four labels are shared across every sheet in the workbook (one per
# year, 2000-2100), stored as shared strings, so the same rename has to be
# applied uniformly on each sheet.

$wb = $excel.ActiveWorkbook

$renames = @{
    "Nd" = "Neodymium"
    "Dy" = "Dysprosium"
    "Cu" = "Copper ores and concentrates"
    "Si" = "Raw silicon"
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    foreach ($cellRef in @("D3", "E3", "F3", "G3")) {
        $cell = $ws.Range($cellRef)
        $current = $cell.Value()
        if ($renames.ContainsKey($current)) {
            $cell.Value = $renames[$current]
        }
    }
}
